# Updated cryptos list - applies new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.417.14"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.922.58"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "376.84"
$ws.Range("E5").Value = "  +6.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.90"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.02"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.31"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").Value = "3.385.15"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.34"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "2.918.54"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.927"
$ws.Range("E17").Value = "  -8.74%  "
$ws.Range("D18").Value = "51.328.36"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.42"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("E27").Value = "  -4.55%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.65"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  +6.87%  "
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.80"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.33"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "34.00"
$ws.Range("E36").Value = "  -5.43%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  -9.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.92"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("E41").Value = "  -10.61%  "
$ws.Range("E42").Value = "  -7.88%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.71"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.53"
$ws.Range("E45").Value = "  -5.59%  "
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.270"
$ws.Range("E47").Value = "  +11.05%  "
$ws.Range("D48").Value = "2.024.66"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.15"
$ws.Range("E50").Value = "  -5.06%  "
$ws.Range("D51").Value = "3.204.37"
$ws.Range("E51").Value = "  -2.80%  "
